# LDLC price tracker: a new timestamped scrape column is inserted right
# before the trailing "nom" / "url_produit" metadata columns, shifting
# those two columns one position to the right. The freshly inserted
# column re-uses the most-recent previously recorded price (the old last
# column, now immediately to its left) as its value for every product
# row, mirroring the "no price movement since last scrape" rows already
# present throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at DO. Excel shifts the old DO -> DP and DP -> DQ,
# and grows the sheet dimension (DP206 -> DQ206) automatically.
$ws.Columns("DO:DO").Insert()

# Header row: label the freshly inserted column with its scrape timestamp.
$ws.Range("DO1").Value = "2026-02-02 05:56:39"

# Data rows 2-80 previously ended with a populated price in column DN
# (now still DN, untouched). Propagate that last known price into the
# brand-new DO column for each of those rows.
$lastPricedRow = 80
for ($r = 2; $r -le $lastPricedRow; $r++) {
    $prevPrice = $ws.Range("DN" + $r).Value()
    $ws.Range("DO" + $r).Value = $prevPrice
}

# Rows 81-206 have no recorded price yet in DN, so the new DO cell stays
# blank as well (matches the rest of the still-empty history columns).

# Two product names recorded slightly differently between scrapes
# (price suffix appended/removed by the source site) - align the shifted
# "nom" column (now DP) with the freshly scraped text for those rows.
$ws.Range("DP29").Value = "Apple iPhone 17 256 Go Lavande"
$ws.Range("DP70").Value = "Apple iPhone 17 Pro Max 512 Go Bleu Intense1 729€00"
